# backup before dimension reduction
# Decrement the numeric suffix of each "qN" label in column A (rows 2-97)
# so that q1 -> q0, q2 -> q1, ..., q96 -> q95.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $txt = [string]$cell.Value2
    $n = [int]($txt.Substring(1))
    $cell.Value = "q" + ($n - 1)
}
